$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Id (A), Ost (Q), Nord (R) values between row 2 and row 3
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$q2 = $ws.Range("Q2").Value2
$q3 = $ws.Range("Q3").Value2
$r2 = $ws.Range("R2").Value2
$r3 = $ws.Range("R3").Value2

$ws.Range("A2").Value2 = $a3
$ws.Range("A3").Value2 = $a2

$ws.Range("Q2").Value2 = $q3
$ws.Range("Q3").Value2 = $q2

$ws.Range("R2").Value2 = $r3
$ws.Range("R3").Value2 = $r2
